# Auto-generated edit script applying the Seraph_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific leve rows
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = ""

$ws.Range("H92").Value = 428.63635
$ws.Range("I92").Value = 421.25
$ws.Range("K92").Value = 421.25
$ws.Range("M92").Value = 826.75

$ws.Range("H113").Value = 4000
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""

$ws.Range("H116").Value = 2927.4285
$ws.Range("I116").Value = 2581
$ws.Range("K116").Value = 2581
$ws.Range("M116").Value = 861

$ws.Range("H138").Value = 3788.9387
$ws.Range("I138").Value = 2457.2222
$ws.Range("J138").Value = 5423.3184
$ws.Range("K138").Value = 7371.6666
$ws.Range("L138").Value = 16269.9552
$ws.Range("M138").Value = -2231.6666
$ws.Range("N138").Value = -26549.9552

$ws.Range("H141").Value = 1465.4
$ws.Range("I141").Value = 1424.375
$ws.Range("K141").Value = 4273.125
$ws.Range("M141").Value = 906.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 207153.73
$ws.Range("I32").Value = 2251.1707
$ws.Range("K32").Value = 2251.1707
$ws.Range("M32").Value = -1964.1707

$ws.Range("H35").Value = 8000
$ws.Range("J35").Value = 8000
$ws.Range("L35").Value = 8000
$ws.Range("N35").Value = -8812

$ws.Range("H122").Value = 4531.636
$ws.Range("I122").Value = 3683.3333
$ws.Range("K122").Value = 11049.9999
$ws.Range("M122").Value = -8599.999899999999

$ws.Range("H132").Value = 1865.4
$ws.Range("I132").Value = 1865.4
$ws.Range("K132").Value = 5596.200000000001
$ws.Range("M132").Value = -3066.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 898.7
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""

$ws.Range("H107").Value = 2292.1
$ws.Range("I107").Value = 1988.7142
$ws.Range("K107").Value = 1988.7142
$ws.Range("M107").Value = -68.71419999999989

$ws.Range("H134").Value = 3457.889
$ws.Range("I134").Value = 2858.8572
$ws.Range("K134").Value = 8576.571599999999
$ws.Range("M134").Value = -6041.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 631.25
$ws.Range("I33").Value = 631.25
$ws.Range("K33").Value = 631.25
$ws.Range("M33").Value = -252.25

$ws.Range("H41").Value = 10119.6
$ws.Range("I41").Value = 6866.3335
$ws.Range("J41").Value = 14999.5
$ws.Range("K41").Value = 6866.3335
$ws.Range("L41").Value = 14999.5
$ws.Range("M41").Value = -6438.3335
$ws.Range("N41").Value = -15855.5

$ws.Range("H134").Value = 2792.64
$ws.Range("I134").Value = 2763.8635
$ws.Range("K134").Value = 8291.5905
$ws.Range("M134").Value = -5756.5905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1275
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 1275
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = ""
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = -4413

$ws.Range("H55").Value = 7787
$ws.Range("J55").Value = 10299.333
$ws.Range("L55").Value = 30897.999
$ws.Range("N55").Value = -31251.999

$ws.Range("H68").Value = 1773.5938
$ws.Range("I68").Value = 2140
$ws.Range("J68").Value = 1705.7407
$ws.Range("K68").Value = 6420
$ws.Range("L68").Value = 5117.2221
$ws.Range("M68").Value = -5609
$ws.Range("N68").Value = -6739.2221

$ws.Range("H71").Value = 1773.5938
$ws.Range("I71").Value = 2140
$ws.Range("J71").Value = 1705.7407
$ws.Range("K71").Value = 19260
$ws.Range("L71").Value = 15351.6663
$ws.Range("M71").Value = -15204
$ws.Range("N71").Value = -23463.6663

$ws.Range("H109").Value = 1361.3572
$ws.Range("I109").Value = 845.125
$ws.Range("K109").Value = 2535.375
$ws.Range("M109").Value = -1495.375

$ws.Range("H132").Value = 3707.6924
$ws.Range("J132").Value = 6359.2
$ws.Range("L132").Value = 57232.8
$ws.Range("N132").Value = -62292.8

$ws.Range("H133").Value = 5000
$ws.Range("J133").Value = 5000
$ws.Range("L133").Value = 15000
$ws.Range("N133").Value = -25120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2542.182
$ws.Range("J97").Value = 7740
$ws.Range("L97").Value = 7740
$ws.Range("N97").Value = -8732

$ws.Range("H102").Value = 4149.4375
$ws.Range("I102").Value = 3374.25
$ws.Range("K102").Value = 3374.25
$ws.Range("M102").Value = -1752.25

$ws.Range("H107").Value = 1734.3846
$ws.Range("I107").Value = 1542.8572
$ws.Range("J107").Value = 1957.8334
$ws.Range("K107").Value = 1542.8572
$ws.Range("L107").Value = 1957.8334
$ws.Range("M107").Value = 377.1428000000001
$ws.Range("N107").Value = -5797.8334

$ws.Range("H122").Value = 64871.688
$ws.Range("I122").Value = 2529.8667
$ws.Range("K122").Value = 7589.6001
$ws.Range("M122").Value = -5139.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 21749.5
$ws.Range("I16").Value = 21749.5
$ws.Range("K16").Value = 21749.5
$ws.Range("M16").Value = -21579.5

$ws.Range("H22").Value = 914.6
$ws.Range("I22").Value = 914.6
$ws.Range("K22").Value = 914.6
$ws.Range("M22").Value = -619.6

$ws.Range("H27").Value = 914.6
$ws.Range("I27").Value = 914.6
$ws.Range("K27").Value = 914.6
$ws.Range("M27").Value = -807.6

$ws.Range("H136").Value = 5265.25
$ws.Range("I136").Value = 4749.4546
$ws.Range("J136").Value = 6400
$ws.Range("K136").Value = 14248.3638
$ws.Range("L136").Value = 19200
$ws.Range("M136").Value = -11698.3638
$ws.Range("N136").Value = -24300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 102000
$ws.Range("J98").Value = 102000
$ws.Range("L98").Value = 102000
$ws.Range("N98").Value = -107990

$ws.Range("H132").Value = 794.087
$ws.Range("I132").Value = 698.8823
$ws.Range("J132").Value = 1063.8334
$ws.Range("K132").Value = 2096.6469
$ws.Range("L132").Value = 3191.5002
$ws.Range("M132").Value = 433.3531000000003
$ws.Range("N132").Value = -8251.5002

$ws.Range("H136").Value = 7143.3
$ws.Range("I136").Value = 6603.6665
$ws.Range("K136").Value = 19810.9995
$ws.Range("M136").Value = -17260.9995
